$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A98").Value = "StudentOpticalFormGradeMismatch"
$ws.Range("B98").Value = "En son okunan optik form  ile ondan once okunan optik form farkli siniflara ait. Bu gibi durumlar raporlardaki ortalama ve siralamayi etkileyebilir."

# Match the formatting used by other rows: column A style mirrors A97 (s="2"),
# column B style mirrors B90 (s="1").
$ws.Range("A97").Copy()
$ws.Range("A98").PasteSpecial(-4122)
$ws.Range("B90").Copy()
$ws.Range("B98").PasteSpecial(-4122)
$excel.CutCopyMode = $false
